$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 3.75
$ws.Range("T2").Value = 1.25
